$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2993.3333
$ws.Range("J51").Value = 2993.3333
$ws.Range("L51").Value = 2993.3333
$ws.Range("N51").Value = -3961.3333

$ws.Range("H111").Value = 10000965
$ws.Range("I111").Value = 16667120
$ws.Range("K111").Value = 50001360
$ws.Range("M111").Value = -49998293

$ws.Range("H132").Value = 580.5185
$ws.Range("I132").Value = 520.1711
$ws.Range("J132").Value = 1497.8
$ws.Range("K132").Value = 1560.5133
$ws.Range("L132").Value = 4493.4
$ws.Range("M132").Value = 969.4866999999999
$ws.Range("N132").Value = -9553.4

$ws.Range("H141").Value = 3312.2
$ws.Range("I141").Value = 1199.6666
$ws.Range("J141").Value = 6481
$ws.Range("K141").Value = 3598.9998
$ws.Range("L141").Value = 19443
$ws.Range("M141").Value = 1581.0002
$ws.Range("N141").Value = -29803

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2856.1619
$ws.Range("I32").Value = 2190.7812
$ws.Range("J32").Value = 13502.25
$ws.Range("K32").Value = 2190.7812
$ws.Range("L32").Value = 13502.25
$ws.Range("M32").Value = -1903.7812
$ws.Range("N32").Value = -14076.25

$ws.Range("H74").Value = 2418.25
$ws.Range("I74").Value = 1599.5
$ws.Range("J74").Value = 2691.1667
$ws.Range("K74").Value = 1599.5
$ws.Range("L74").Value = 2691.1667
$ws.Range("M74").Value = -725.5
$ws.Range("N74").Value = -4439.1667

$ws.Range("H77").Value = 2418.25
$ws.Range("I77").Value = 1599.5
$ws.Range("J77").Value = 2691.1667
$ws.Range("K77").Value = 7997.5
$ws.Range("L77").Value = 13455.8335
$ws.Range("M77").Value = -3629.5
$ws.Range("N77").Value = -22191.8335

$ws.Range("H110").Value = 1949.76
$ws.Range("I110").Value = 1634.8572
$ws.Range("J110").Value = 3603
$ws.Range("K110").Value = 1634.8572
$ws.Range("L110").Value = 3603
$ws.Range("M110").Value = 410.1428000000001
$ws.Range("N110").Value = -7693

$ws.Range("H122").Value = 1136.3334
$ws.Range("I122").Value = 1297.7858
$ws.Range("J122").Value = 813.4286
$ws.Range("K122").Value = 3893.3574
$ws.Range("L122").Value = 2440.2858
$ws.Range("M122").Value = -1443.3574
$ws.Range("N122").Value = -7340.2858

$ws.Range("H132").Value = 1929.5714
$ws.Range("I132").Value = 1377.125
$ws.Range("J132").Value = 3697.4
$ws.Range("K132").Value = 4131.375
$ws.Range("L132").Value = 11092.2
$ws.Range("M132").Value = -1601.375
$ws.Range("N132").Value = -16152.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2025.7142
$ws.Range("I99").Value = 1857.5
$ws.Range("J99").Value = 2250
$ws.Range("K99").Value = 1857.5
$ws.Range("L99").Value = 2250
$ws.Range("M99").Value = -359.5
$ws.Range("N99").Value = -5246

$ws.Range("H122").Value = 48666.668
$ws.Range("J122").Value = 48666.668
$ws.Range("L122").Value = 48666.668
$ws.Range("N122").Value = -58466.668

$ws.Range("H134").Value = 7014.7856
$ws.Range("I134").Value = 7631.2607
$ws.Range("J134").Value = 4179
$ws.Range("K134").Value = 22893.7821
$ws.Range("L134").Value = 12537
$ws.Range("M134").Value = -20358.7821
$ws.Range("N134").Value = -17607

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1090.7693
$ws.Range("J22").Value = 1327.9
$ws.Range("L22").Value = 1327.9
$ws.Range("N22").Value = -2027.9

$ws.Range("H31").Value = 1886
$ws.Range("I31").Value = 1748.4
$ws.Range("K31").Value = 1748.4
$ws.Range("M31").Value = -1453.4

$ws.Range("H34").Value = 1886
$ws.Range("I34").Value = 1748.4
$ws.Range("K34").Value = 1748.4
$ws.Range("M34").Value = -1546.4

$ws.Range("H58").Value = 1176972.2
$ws.Range("I58").Value = 1611668.4
$ws.Range("J58").Value = 3292.5
$ws.Range("K58").Value = 1611668.4
$ws.Range("L58").Value = 3292.5
$ws.Range("M58").Value = -1611465.4
$ws.Range("N58").Value = -3698.5

$ws.Range("H99").Value = 1867.591
$ws.Range("I99").Value = 1891.5385
$ws.Range("J99").Value = 1833
$ws.Range("K99").Value = 1891.5385
$ws.Range("L99").Value = 1833
$ws.Range("M99").Value = -393.5385000000001
$ws.Range("N99").Value = -4829

$ws.Range("H112").Value = 25000
$ws.Range("J112").Value = 25000
$ws.Range("L112").Value = 25000
$ws.Range("N112").Value = -27954

$ws.Range("H126").Value = 1867.591
$ws.Range("I126").Value = 1891.5385
$ws.Range("J126").Value = 1833
$ws.Range("K126").Value = 5674.6155
$ws.Range("L126").Value = 5499
$ws.Range("M126").Value = -3204.6155
$ws.Range("N126").Value = -10439

$ws.Range("H132").Value = 2002.1892
$ws.Range("I132").Value = 1078.5834
$ws.Range("K132").Value = 3235.7502
$ws.Range("M132").Value = -705.7501999999999

$ws.Range("H134").Value = 2195.88
$ws.Range("I134").Value = 1886.8695
$ws.Range("J134").Value = 5749.5
$ws.Range("K134").Value = 5660.6085
$ws.Range("L134").Value = 17248.5
$ws.Range("M134").Value = -3125.6085
$ws.Range("N134").Value = -22318.5

$ws.Range("H136").Value = 1176972.2
$ws.Range("I136").Value = 1611668.4
$ws.Range("J136").Value = 3292.5
$ws.Range("K136").Value = 4835005.199999999
$ws.Range("L136").Value = 9877.5
$ws.Range("M136").Value = -4832455.199999999
$ws.Range("N136").Value = -14977.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 399.44446
$ws.Range("J98").Value = 491.33334
$ws.Range("L98").Value = 1474.00002
$ws.Range("N98").Value = -4470.000019999999

$ws.Range("H104").Value = 3823.9412
$ws.Range("J104").Value = 4325.5
$ws.Range("L104").Value = 12976.5
$ws.Range("N104").Value = -18218.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 21926.4
$ws.Range("J92").Value = 21926.4
$ws.Range("L92").Value = 21926.4
$ws.Range("N92").Value = -25670.4

$ws.Range("H97").Value = 2272.8
$ws.Range("I97").Value = 2103.625
$ws.Range("J97").Value = 2949.5
$ws.Range("K97").Value = 2103.625
$ws.Range("L97").Value = 2949.5
$ws.Range("M97").Value = -1607.625
$ws.Range("N97").Value = -3941.5

$ws.Range("H126").Value = 2461357
$ws.Range("I126").Value = 6175466.5
$ws.Range("K126").Value = 18526399.5
$ws.Range("M126").Value = -18523929.5

$ws.Range("H132").Value = 2140278.5
$ws.Range("I132").Value = 3848533.5
$ws.Range("K132").Value = 11545600.5
$ws.Range("M132").Value = -11543070.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3433.2
$ws.Range("I7").Value = 2703.5557
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 2703.5557
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -2591.5557
$ws.Range("N7").Value = -10224

$ws.Range("H22").Value = 2233.7778
$ws.Range("I22").Value = 2507
$ws.Range("J22").Value = 2059.9092
$ws.Range("K22").Value = 2507
$ws.Range("L22").Value = 2059.9092
$ws.Range("M22").Value = -2212
$ws.Range("N22").Value = -2649.9092

$ws.Range("H27").Value = 2233.7778
$ws.Range("I27").Value = 2507
$ws.Range("J27").Value = 2059.9092
$ws.Range("K27").Value = 2507
$ws.Range("L27").Value = 2059.9092
$ws.Range("M27").Value = -2400
$ws.Range("N27").Value = -2273.9092

$ws.Range("H40").Value = 3582.3125
$ws.Range("I40").Value = 2281.4
$ws.Range("K40").Value = 2281.4
$ws.Range("M40").Value = -2145.4

$ws.Range("H82").Value = 1557.2858
$ws.Range("I82").Value = 1400.2
$ws.Range("K82").Value = 1400.2
$ws.Range("M82").Value = -1039.2

$ws.Range("H85").Value = 1557.2858
$ws.Range("I85").Value = 1400.2
$ws.Range("K85").Value = 1400.2
$ws.Range("M85").Value = -152.2

$ws.Range("H126").Value = 3433.2
$ws.Range("I126").Value = 2703.5557
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 8110.6671
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -5640.6671
$ws.Range("N126").Value = -34940

$ws.Range("H132").Value = 1982.08
$ws.Range("I132").Value = 1242.4
$ws.Range("K132").Value = 3727.2
$ws.Range("M132").Value = -1197.2

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 24137.8
$ws.Range("J119").Value = 24137.8
$ws.Range("L119").Value = 24137.8
$ws.Range("N119").Value = -33813.8
